$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the new data range so numeric- and date-looking
# strings (dates, times, measurements) are preserved verbatim as text,
# matching the source data which stores every value as a string.
$ws.Range("A510:K536").NumberFormat = "@"

# Row 510
$ws.Range("A510").Value = "stampfenbachstrasse"
$ws.Range("B510").Value = "14/04/2021"
$ws.Range("C510").Value = "19:27:42"
$ws.Range("D510").Value = "Broken clouds"
$ws.Range("E510").Value = "5"
$ws.Range("F510").Value = "45"
$ws.Range("G510").Value = "16.7"
$ws.Range("H510").Value = "1024"
$ws.Range("I510").Value = "N/A"
$ws.Range("J510").Value = "10.2"
$ws.Range("K510").Value = " 30 "

# Row 511
$ws.Range("A511").Value = "stampfenbachstrasse"
$ws.Range("B511").Value = "14/04/2021"
$ws.Range("C511").Value = "20:00:03"
$ws.Range("D511").Value = "Broken clouds"
$ws.Range("E511").Value = "5"
$ws.Range("F511").Value = "45"
$ws.Range("G511").Value = "16.7"
$ws.Range("H511").Value = "1024"
$ws.Range("I511").Value = "N/A"
$ws.Range("J511").Value = "10.2"
$ws.Range("K511").Value = " 30 "

# Row 512
$ws.Range("A512").Value = "stampfenbachstrasse"
$ws.Range("B512").Value = "14/04/2021"
$ws.Range("C512").Value = "21:00:02"
$ws.Range("D512").Value = "Broken clouds"
$ws.Range("E512").Value = "5"
$ws.Range("F512").Value = "45"
$ws.Range("G512").Value = "16.7"
$ws.Range("H512").Value = "1024"
$ws.Range("I512").Value = "6"
$ws.Range("J512").Value = "11.3"
$ws.Range("K512").Value = " 25 "

# Row 513
$ws.Range("A513").Value = "stampfenbachstrasse"
$ws.Range("B513").Value = "14/04/2021"
$ws.Range("C513").Value = "21:30:03"
$ws.Range("D513").Value = "Scattered clouds"
$ws.Range("E513").Value = "3"
$ws.Range("F513").Value = "65"
$ws.Range("G513").Value = "7.4"
$ws.Range("H513").Value = "1025"
$ws.Range("I513").Value = "N/A"
$ws.Range("J513").Value = "11.5"
$ws.Range("K513").Value = " 26 "

# Row 514
$ws.Range("A514").Value = "stampfenbachstrasse"
$ws.Range("B514").Value = "14/04/2021"
$ws.Range("C514").Value = "22:00:03"
$ws.Range("D514").Value = "Scattered clouds"
$ws.Range("E514").Value = "3"
$ws.Range("F514").Value = "65"
$ws.Range("G514").Value = "7.4"
$ws.Range("H514").Value = "1025"
$ws.Range("I514").Value = "N/A"
$ws.Range("J514").Value = "11.5"
$ws.Range("K514").Value = " 26 "

# Row 515
$ws.Range("A515").Value = "stampfenbachstrasse"
$ws.Range("B515").Value = "14/04/2021"
$ws.Range("C515").Value = "22:30:03"
$ws.Range("D515").Value = "Scattered clouds"
$ws.Range("E515").Value = "3"
$ws.Range("F515").Value = "65"
$ws.Range("G515").Value = "7.4"
$ws.Range("H515").Value = "1025"
$ws.Range("I515").Value = "8"
$ws.Range("J515").Value = "11.4"
$ws.Range("K515").Value = " 33 "

# Row 516
$ws.Range("A516").Value = "stampfenbachstrasse"
$ws.Range("B516").Value = "14/04/2021"
$ws.Range("C516").Value = "23:00:03"
$ws.Range("D516").Value = "Scattered clouds"
$ws.Range("E516").Value = "3"
$ws.Range("F516").Value = "65"
$ws.Range("G516").Value = "7.4"
$ws.Range("H516").Value = "1025"
$ws.Range("I516").Value = "8"
$ws.Range("J516").Value = "11.4"
$ws.Range("K516").Value = " 33 "

# Row 517
$ws.Range("A517").Value = "stampfenbachstrasse"
$ws.Range("B517").Value = "14/04/2021"
$ws.Range("C517").Value = "23:30:03"
$ws.Range("D517").Value = "Broken clouds"
$ws.Range("E517").Value = "2"
$ws.Range("F517").Value = "65"
$ws.Range("G517").Value = "3.9"
$ws.Range("H517").Value = "1025"
$ws.Range("I517").Value = "N/A"
$ws.Range("J517").Value = "10.8"
$ws.Range("K517").Value = " 28 "

# Row 518
$ws.Range("A518").Value = "stampfenbachstrasse"
$ws.Range("B518").Value = "15/04/2021"
$ws.Range("C518").Value = "00:00:03"
$ws.Range("D518").Value = "Broken clouds"
$ws.Range("E518").Value = "2"
$ws.Range("F518").Value = "65"
$ws.Range("G518").Value = "3.9"
$ws.Range("H518").Value = "1025"
$ws.Range("I518").Value = "N/A"
$ws.Range("J518").Value = "10.8"
$ws.Range("K518").Value = " 28 "

# Row 519
$ws.Range("A519").Value = "stampfenbachstrasse"
$ws.Range("B519").Value = "15/04/2021"
$ws.Range("C519").Value = "00:30:03"
$ws.Range("D519").Value = "Broken clouds"
$ws.Range("E519").Value = "2"
$ws.Range("F519").Value = "65"
$ws.Range("G519").Value = "3.9"
$ws.Range("H519").Value = "1025"
$ws.Range("I519").Value = "N/A"
$ws.Range("J519").Value = "10.8"
$ws.Range("K519").Value = " 28 "

# Row 520
$ws.Range("A520").Value = "stampfenbachstrasse"
$ws.Range("B520").Value = "15/04/2021"
$ws.Range("C520").Value = "01:00:03"
$ws.Range("D520").Value = "Broken clouds"
$ws.Range("E520").Value = "2"
$ws.Range("F520").Value = "65"
$ws.Range("G520").Value = "3.9"
$ws.Range("H520").Value = "1025"
$ws.Range("I520").Value = "N/A"
$ws.Range("J520").Value = "10.8"
$ws.Range("K520").Value = " 28 "

# Row 521
$ws.Range("A521").Value = "stampfenbachstrasse"
$ws.Range("B521").Value = "15/04/2021"
$ws.Range("C521").Value = "01:30:03"
$ws.Range("D521").Value = "Broken clouds"
$ws.Range("E521").Value = "2"
$ws.Range("F521").Value = "65"
$ws.Range("G521").Value = "3.9"
$ws.Range("H521").Value = "1025"
$ws.Range("I521").Value = "N/A"
$ws.Range("J521").Value = "10.8"
$ws.Range("K521").Value = " 28 "

# Row 522
$ws.Range("A522").Value = "stampfenbachstrasse"
$ws.Range("B522").Value = "15/04/2021"
$ws.Range("C522").Value = "02:00:03"
$ws.Range("D522").Value = "Broken clouds"
$ws.Range("E522").Value = "2"
$ws.Range("F522").Value = "65"
$ws.Range("G522").Value = "3.9"
$ws.Range("H522").Value = "1025"
$ws.Range("I522").Value = "N/A"
$ws.Range("J522").Value = "10.8"
$ws.Range("K522").Value = " 28 "

# Row 523
$ws.Range("A523").Value = "stampfenbachstrasse"
$ws.Range("B523").Value = "15/04/2021"
$ws.Range("C523").Value = "02:30:03"
$ws.Range("D523").Value = "Scattered clouds"
$ws.Range("E523").Value = "0"
$ws.Range("F523").Value = "69"
$ws.Range("G523").Value = "1.2"
$ws.Range("H523").Value = "1025"
$ws.Range("I523").Value = "9"
$ws.Range("J523").Value = "12.6"
$ws.Range("K523").Value = " 37 "

# Row 524
$ws.Range("A524").Value = "stampfenbachstrasse"
$ws.Range("B524").Value = "15/04/2021"
$ws.Range("C524").Value = "03:00:03"
$ws.Range("D524").Value = "Scattered clouds"
$ws.Range("E524").Value = "0"
$ws.Range("F524").Value = "69"
$ws.Range("G524").Value = "1.2"
$ws.Range("H524").Value = "1025"
$ws.Range("I524").Value = "9"
$ws.Range("J524").Value = "12.6"
$ws.Range("K524").Value = " 37 "

# Row 525
$ws.Range("A525").Value = "stampfenbachstrasse"
$ws.Range("B525").Value = "15/04/2021"
$ws.Range("C525").Value = "03:30:03"
$ws.Range("D525").Value = "Scattered clouds"
$ws.Range("E525").Value = "0"
$ws.Range("F525").Value = "69"
$ws.Range("G525").Value = "1.2"
$ws.Range("H525").Value = "1025"
$ws.Range("I525").Value = "9"
$ws.Range("J525").Value = "N/A"
$ws.Range("K525").Value = " 37 "

# Row 526
$ws.Range("A526").Value = "stampfenbachstrasse"
$ws.Range("B526").Value = "15/04/2021"
$ws.Range("C526").Value = "04:00:03"
$ws.Range("D526").Value = "Scattered clouds"
$ws.Range("E526").Value = "0"
$ws.Range("F526").Value = "69"
$ws.Range("G526").Value = "1.2"
$ws.Range("H526").Value = "1025"
$ws.Range("I526").Value = "9"
$ws.Range("J526").Value = "N/A"
$ws.Range("K526").Value = " 37 "

# Row 527
$ws.Range("A527").Value = "stampfenbachstrasse"
$ws.Range("B527").Value = "15/04/2021"
$ws.Range("C527").Value = "04:30:03"
$ws.Range("D527").Value = "Scattered clouds"
$ws.Range("E527").Value = "0"
$ws.Range("F527").Value = "69"
$ws.Range("G527").Value = "1.2"
$ws.Range("H527").Value = "1025"
$ws.Range("I527").Value = "N/A"
$ws.Range("J527").Value = "13.2"
$ws.Range("K527").Value = " 22 "

# Row 528
$ws.Range("A528").Value = "stampfenbachstrasse"
$ws.Range("B528").Value = "15/04/2021"
$ws.Range("C528").Value = "05:00:03"
$ws.Range("D528").Value = "Scattered clouds"
$ws.Range("E528").Value = "0"
$ws.Range("F528").Value = "69"
$ws.Range("G528").Value = "1.2"
$ws.Range("H528").Value = "1025"
$ws.Range("I528").Value = "N/A"
$ws.Range("J528").Value = "13.2"
$ws.Range("K528").Value = " 22 "

# Row 529
$ws.Range("A529").Value = "stampfenbachstrasse"
$ws.Range("B529").Value = "15/04/2021"
$ws.Range("C529").Value = "05:30:03"
$ws.Range("D529").Value = "Scattered clouds"
$ws.Range("E529").Value = "-1"
$ws.Range("F529").Value = "93"
$ws.Range("G529").Value = "1.6"
$ws.Range("H529").Value = "1024"
$ws.Range("I529").Value = "9"
$ws.Range("J529").Value = "13.4"
$ws.Range("K529").Value = " 37 "

# Row 530
$ws.Range("A530").Value = "stampfenbachstrasse"
$ws.Range("B530").Value = "15/04/2021"
$ws.Range("C530").Value = "06:00:03"
$ws.Range("D530").Value = "Scattered clouds"
$ws.Range("E530").Value = "-1"
$ws.Range("F530").Value = "93"
$ws.Range("G530").Value = "1.6"
$ws.Range("H530").Value = "1024"
$ws.Range("I530").Value = "9"
$ws.Range("J530").Value = "13.4"
$ws.Range("K530").Value = " 37 "

# Row 531
$ws.Range("A531").Value = "stampfenbachstrasse"
$ws.Range("B531").Value = "15/04/2021"
$ws.Range("C531").Value = "06:30:03"
$ws.Range("D531").Value = "Scattered clouds"
$ws.Range("E531").Value = "-1"
$ws.Range("F531").Value = "93"
$ws.Range("G531").Value = "1.6"
$ws.Range("H531").Value = "1024"
$ws.Range("I531").Value = "10"
$ws.Range("J531").Value = "14"
$ws.Range("K531").Value = " 41 "

# Row 532
$ws.Range("A532").Value = "stampfenbachstrasse"
$ws.Range("B532").Value = "15/04/2021"
$ws.Range("C532").Value = "07:00:03"
$ws.Range("D532").Value = "Scattered clouds"
$ws.Range("E532").Value = "-1"
$ws.Range("F532").Value = "93"
$ws.Range("G532").Value = "1.6"
$ws.Range("H532").Value = "1024"
$ws.Range("I532").Value = "10"
$ws.Range("J532").Value = "14"
$ws.Range("K532").Value = " 41 "

# Row 533
$ws.Range("A533").Value = "stampfenbachstrasse"
$ws.Range("B533").Value = "15/04/2021"
$ws.Range("C533").Value = "07:30:03"
$ws.Range("D533").Value = "Few clouds"
$ws.Range("E533").Value = "-1"
$ws.Range("F533").Value = "93"
$ws.Range("G533").Value = "2.1"
$ws.Range("H533").Value = "1024"
$ws.Range("I533").Value = "N/A"
$ws.Range("J533").Value = "15.8"
$ws.Range("K533").Value = " 14 "

# Row 534
$ws.Range("A534").Value = "stampfenbachstrasse"
$ws.Range("B534").Value = "15/04/2021"
$ws.Range("C534").Value = "08:00:03"
$ws.Range("D534").Value = "Few clouds"
$ws.Range("E534").Value = "-1"
$ws.Range("F534").Value = "93"
$ws.Range("G534").Value = "2.1"
$ws.Range("H534").Value = "1024"
$ws.Range("I534").Value = "N/A"
$ws.Range("J534").Value = "15.8"
$ws.Range("K534").Value = " 14 "

# Row 535
$ws.Range("A535").Value = "stampfenbachstrasse"
$ws.Range("B535").Value = "15/04/2021"
$ws.Range("C535").Value = "08:30:03"
$ws.Range("D535").Value = "Few clouds"
$ws.Range("E535").Value = "-1"
$ws.Range("F535").Value = "93"
$ws.Range("G535").Value = "2.1"
$ws.Range("H535").Value = "1024"
$ws.Range("I535").Value = "11"
$ws.Range("J535").Value = "16.9"
$ws.Range("K535").Value = " 45 "

# Row 536
$ws.Range("A536").Value = "stampfenbachstrasse"
$ws.Range("B536").Value = "15/04/2021"
$ws.Range("C536").Value = "09:00:03"
$ws.Range("D536").Value = "Few clouds"
$ws.Range("E536").Value = "-1"
$ws.Range("F536").Value = "93"
$ws.Range("G536").Value = "2.1"
$ws.Range("H536").Value = "1024"
$ws.Range("I536").Value = "11"
$ws.Range("J536").Value = "16.9"
$ws.Range("K536").Value = " 45 "
